$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.091.48"
$ws.Range("E2").Value = "  +0.20%  "

$ws.Range("D3").Value = "1.834.12"
$ws.Range("E3").Value = "  +0.03%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9992"
$ws.Range("E4").Value = "  +0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "243.66"
$ws.Range("E5").Value = "  +0.55%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6280"
$ws.Range("E6").Value = "  +0.08%  "

$ws.Range("E7").Value = "  +0.17%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07520"

$ws.Range("E9").Value = "  -0.21%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.23"
$ws.Range("E10").Value = "  +3.02%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07691"
$ws.Range("E11").Value = "  -0.29%  "

$ws.Range("D12").Value = "1.832.89"
$ws.Range("E12").Value = "  -0.39%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.003"
$ws.Range("E13").Value = "  +0.93%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6679"
$ws.Range("E14").Value = "  +0.36%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "82.74"
$ws.Range("E15").Value = "  -0.06%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.000009394"
$ws.Range("E16").Value = "  -7.34%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.986"
$ws.Range("E17").Value = "  -1.02%  "

$ws.Range("D18").Value = "29.095.77"

$ws.Range("D19").Value = "2.079.71"
$ws.Range("E19").Value = "  -0.26%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.59"
$ws.Range("E20").Value = "  +1.88%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "223.54"
$ws.Range("E21").Value = "  -1.35%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.004"
$ws.Range("E22").Value = "  +0.50%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.101"
$ws.Range("E23").Value = "  -1.11%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.000"
$ws.Range("E24").Value = "  +0.10%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "159.65"
$ws.Range("E25").Value = "  +0.85%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1399"
$ws.Range("E26").Value = "  +1.70%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.495"
$ws.Range("E27").Value = "  -0.07%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.93"
$ws.Range("E28").Value = "  +0.14%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.500"
$ws.Range("E29").Value = "  +0.69%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05666"
$ws.Range("E30").Value = "  +8.37%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.156"
$ws.Range("E31").Value = "  +1.07%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.080"
$ws.Range("E32").Value = "  +1.55%  "

$ws.Range("E33").Value = "  +1.13%  "

$ws.Range("B34").Value = "LidoDAOToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.841"
$ws.Range("E34").Value = "  -0.16%  "

$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7426"
$ws.Range("E35").Value = "  +0.88%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.140"
$ws.Range("E36").Value = "  +0.03%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.671"
$ws.Range("E37").Value = "  -1.35%  "

$ws.Range("E38").Value = "  +0.12%  "

$ws.Range("D39").Value = "1.219.92"
$ws.Range("E39").Value = "  -1.69%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01783"
$ws.Range("E40").Value = "  -0.15%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.531"
$ws.Range("E41").Value = "  +2.85%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8922"
$ws.Range("E42").Value = "  -0.32%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.001"
$ws.Range("E43").Value = "  +0.17%  "

$ws.Range("B44").Value = "BabyDogeCoin"
$ws.Range("C44").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.00000000128"
$ws.Range("E44").Value = "  +4.05%  "

$ws.Range("B45").Value = "Quant"
$ws.Range("C45").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "101.93"
$ws.Range("E45").Value = "  +0.20%  "

$ws.Range("B46").Value = "RocketPoolETH"
$ws.Range("C46").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D46").Value = "1.976.38"
$ws.Range("E46").Value = "  -0.27%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "65.82"
$ws.Range("E47").Value = "  +2.50%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5086"

$ws.Range("B49").Value = "TheSandbox"
$ws.Range("C49").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.4081"
$ws.Range("E49").Value = "  +1.04%  "

$ws.Range("B50").Value = "XinFinNetwork"
$ws.Range("C50").Value = "https://coinranking.com/coin/77jGXSqWJ1ofG+xinfinnetwork-xdc"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07515"
$ws.Range("E50").Value = "  +7.90%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.988"
$ws.Range("E51").Value = "  +0.65%  "
